$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = -12.8801
$ws.Range("A9").Value = -21.8389
$ws.Range("C12").Value = -10.6982
$ws.Range("A13").Value = -22.1798
$ws.Range("C14").Value = -14.0201
$ws.Range("A16").Value = -21.55879999999999
$ws.Range("A18").Value = -22.33850000000002
$ws.Range("C19").Value = -12.01390000000001
$ws.Range("A20").Value = -21.55159999999999
$ws.Range("A26").Value = -21.58119999999998
$ws.Range("C26").Value = -12.9093
$ws.Range("A27").Value = -21.63229999999999
$ws.Range("C27").Value = -12.5743
$ws.Range("A29").Value = -20.63589999999998
$ws.Range("C29").Value = -11.6286
$ws.Range("A35").Value = -21.43799999999997
$ws.Range("A36").Value = -21.12949999999999
$ws.Range("C37").Value = -13.5447
$ws.Range("C38").Value = -12.4358
$ws.Range("A45").Value = -21.57649999999998
$ws.Range("C47").Value = -11.9473
$ws.Range("C51").Value = -12.41939999999999
$ws.Range("C52").Value = -11.35750000000001
$ws.Range("A55").Value = -22.44839999999999
$ws.Range("C55").Value = -13.19979999999999
$ws.Range("A57").Value = -22.26550000000001
$ws.Range("A69").Value = -21.73359999999998
$ws.Range("C69").Value = -11.1957
$ws.Range("C70").Value = -11.93140000000001
$ws.Range("A76").Value = -19.3324
$ws.Range("C76").Value = -12.81100000000001
$ws.Range("A78").Value = -19.96279999999999
$ws.Range("C81").Value = -13.5196
$ws.Range("A82").Value = -22.21100000000001
$ws.Range("A83").Value = -21.9449
$ws.Range("C83").Value = -13.22969999999999
$ws.Range("A93").Value = -20.58619999999998
$ws.Range("C94").Value = -10.2655
$ws.Range("A97").Value = -22.02940000000001
$ws.Range("C100").Value = -12.16
$ws.Range("C102").Value = -13.4775
